$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.532.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.080.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("E11").Value = "  +2.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.389.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.079.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.540.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("E28").Value = "  -3.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  +0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0636"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("E37").Value = "  -2.14%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0231"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0956"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.460.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("E47").Value = "  -0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.272.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
